# Fix errors caused by pandas upgrade
# Extend the cuts_head header range and add a new "survey_code" cut
# (male/female) on the Lookups sheet, shifting the zero_string cell
# two columns to the right (from V1 to X1).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lookups")

# New header cell + value/order pair for the survey_code cut, mirroring
# the pattern used by the existing Gender cut in columns G:H.
$ws.Range("V1").Value = "survey_code"
$ws.Range("V2").Value = "male"
$ws.Range("W2").Value = 1
$ws.Range("V3").Value = "female"
$ws.Range("W3").Value = 2

# The zero_string marker cell moves from V1 to X1 to make room for the
# new survey_code columns.
$ws.Range("X1").Value = 0

# Update the defined names so they reflect the new layout.
$wb.Names("cuts_head").RefersTo = "='Lookups'!`$G`$1:`$W`$1"
$wb.Names("zero_string").RefersTo = "='Lookups'!`$X`$1"
